# Auto-generated cell updates applying the scheduled-runner price refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4167675.8
$ws.Range("I19").Value = 8333756.5
$ws.Range("J19").Value = 1594.6666
$ws.Range("K19").Value = 8333756.5
$ws.Range("L19").Value = 1594.6666
$ws.Range("M19").Value = -8333581.5
$ws.Range("N19").Value = -1944.6666
$ws.Range("H21").Value = 73859.5
$ws.Range("I21").Value = 68000
$ws.Range("J21").Value = 79719
$ws.Range("K21").Value = 68000
$ws.Range("L21").Value = 79719
$ws.Range("M21").Value = -67532
$ws.Range("N21").Value = -80655
$ws.Range("H23").Value = 73859.5
$ws.Range("I23").Value = 68000
$ws.Range("J23").Value = 79719
$ws.Range("K23").Value = 68000
$ws.Range("L23").Value = 79719
$ws.Range("M23").Value = -67766
$ws.Range("N23").Value = -80187
$ws.Range("H28").Value = 18520346
$ws.Range("I28").Value = 29414474
$ws.Range("J28").Value = 327.6
$ws.Range("K28").Value = 29414474
$ws.Range("L28").Value = 327.6
$ws.Range("M28").Value = -29413989
$ws.Range("N28").Value = -1297.6
$ws.Range("H29").Value = 251.5
$ws.Range("I29").Value = 251.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 754.5
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -473.5
$ws.Range("H38").Value = 1575.238
$ws.Range("I38").Value = 142.5
$ws.Range("J38").Value = 3485.5557
$ws.Range("K38").Value = 427.5
$ws.Range("L38").Value = 10456.6671
$ws.Range("M38").Value = -55.5
$ws.Range("N38").Value = -11200.6671
$ws.Range("H58").Value = 1038.5
$ws.Range("I58").Value = 57.75
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 173.25
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -23.25
$ws.Range("N58").Value = -9300
$ws.Range("H64").Value = 2905.2788
$ws.Range("I64").Value = 2920.7778
$ws.Range("J64").Value = 2785.7144
$ws.Range("K64").Value = 2920.7778
$ws.Range("L64").Value = 2785.7144
$ws.Range("M64").Value = -2672.7778
$ws.Range("N64").Value = -3281.7144
$ws.Range("H67").Value = 2905.2788
$ws.Range("I67").Value = 2920.7778
$ws.Range("J67").Value = 2785.7144
$ws.Range("K67").Value = 2920.7778
$ws.Range("L67").Value = 2785.7144
$ws.Range("M67").Value = -2062.7778
$ws.Range("N67").Value = -4501.7144
$ws.Range("H87").Value = 12866.639
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 12866.639
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 12866.639
$ws.Range("N87").Value = -15362.639
$ws.Range("H90").Value = 12866.639
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 12866.639
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 38599.917
$ws.Range("N90").Value = -51079.917
$ws.Range("H106").Value = 90913410
$ws.Range("I106").Value = 333338500
$ws.Range("J106").Value = 3997.625
$ws.Range("K106").Value = 333338500
$ws.Range("L106").Value = 3997.625
$ws.Range("M106").Value = -333337869
$ws.Range("N106").Value = -5259.625
$ws.Range("H113").Value = 7694186
$ws.Range("I113").Value = 10001741
$ws.Range("J113").Value = 2335.3333
$ws.Range("K113").Value = 10001741
$ws.Range("L113").Value = 2335.3333
$ws.Range("M113").Value = -9998487
$ws.Range("N113").Value = -8843.3333
$ws.Range("H132").Value = 3368547.5
$ws.Range("I132").Value = 1595.4193
$ws.Range("J132").Value = 55556304
$ws.Range("K132").Value = 4786.257900000001
$ws.Range("L132").Value = 166668912
$ws.Range("M132").Value = -2256.257900000001
$ws.Range("N132").Value = -166673972
$ws.Range("H138").Value = 3407.4055
$ws.Range("I138").Value = 3969.1875
$ws.Range("J138").Value = 2979.3809
$ws.Range("K138").Value = 11907.5625
$ws.Range("L138").Value = 8938.1427
$ws.Range("M138").Value = -6767.5625
$ws.Range("N138").Value = -19218.1427
$ws.Range("H139").Value = 40404.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 40404.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 40404.332
$ws.Range("N139").Value = -50684.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2872366.8
$ws.Range("I61").Value = 1812692.9
$ws.Range("J61").Value = 4903408.5
$ws.Range("K61").Value = 1812692.9
$ws.Range("L61").Value = 4903408.5
$ws.Range("M61").Value = -1812480.9
$ws.Range("N61").Value = -4903832.5
$ws.Range("H132").Value = 9975369
$ws.Range("I132").Value = 11909157
$ws.Range("J132").Value = 5053000.5
$ws.Range("K132").Value = 35727471
$ws.Range("L132").Value = 15159001.5
$ws.Range("M132").Value = -35724941
$ws.Range("N132").Value = -15164061.5
$ws.Range("H136").Value = 2872366.8
$ws.Range("I136").Value = 1812692.9
$ws.Range("J136").Value = 4903408.5
$ws.Range("K136").Value = 5438078.699999999
$ws.Range("L136").Value = 14710225.5
$ws.Range("M136").Value = -5435528.699999999
$ws.Range("N136").Value = -14715325.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2175.0227
$ws.Range("I68").Value = 655.78125
$ws.Range("J68").Value = 3043.1606
$ws.Range("K68").Value = 1967.34375
$ws.Range("L68").Value = 9129.481800000001
$ws.Range("M68").Value = -1156.34375
$ws.Range("N68").Value = -10751.4818
$ws.Range("H71").Value = 2175.0227
$ws.Range("I71").Value = 655.78125
$ws.Range("J71").Value = 3043.1606
$ws.Range("K71").Value = 5902.03125
$ws.Range("L71").Value = 27388.4454
$ws.Range("M71").Value = -1846.03125
$ws.Range("N71").Value = -35500.4454
$ws.Range("H101").Value = 9999.571
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 9999.571
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 29998.713
$ws.Range("N101").Value = -34866.713
$ws.Range("H131").Value = 39600.348
$ws.Range("I131").Value = 200346
$ws.Range("J131").Value = 1327.5714
$ws.Range("K131").Value = 601038
$ws.Range("L131").Value = 3982.7142
$ws.Range("M131").Value = -595998
$ws.Range("N131").Value = -14062.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6544.6665
$ws.Range("I126").Value = 7920.7144
$ws.Range("J126").Value = 1728.5
$ws.Range("K126").Value = 23762.1432
$ws.Range("L126").Value = 5185.5
$ws.Range("M126").Value = -21292.1432
$ws.Range("N126").Value = -10125.5
$ws.Range("H132").Value = 11393625
$ws.Range("I132").Value = 12746104
$ws.Range("J132").Value = 9094409
$ws.Range("K132").Value = 38238312
$ws.Range("L132").Value = 27283227
$ws.Range("M132").Value = -38235782
$ws.Range("N132").Value = -27288287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H10").Value = 476.6
$ws.Range("I10").Value = 476.6
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 476.6
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -336.6
$ws.Range("H17").Value = 26323.4
$ws.Range("I17").Value = 866.3333
$ws.Range("J17").Value = 64509
$ws.Range("K17").Value = 866.3333
$ws.Range("L17").Value = 64509
$ws.Range("M17").Value = -696.3333
$ws.Range("N17").Value = -64849
$ws.Range("H18").Value = 1247.5
$ws.Range("I18").Value = 995
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 995
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -823
$ws.Range("N18").Value = -1844
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1000
$ws.Range("N30").Value = -1216
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H125").Value = 50711.668
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 50711.668
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 50711.668
$ws.Range("N125").Value = -60551.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 9127.5
$ws.Range("I8").Value = 12166.667
$ws.Range("J8").Value = 10
$ws.Range("K8").Value = 12166.667
$ws.Range("L8").Value = 10
$ws.Range("M8").Value = -12026.667
$ws.Range("N8").Value = -290
$ws.Range("H74").Value = 16256.909
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 16256.909
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 16256.909
$ws.Range("N74").Value = -18128.909
$ws.Range("H77").Value = 16256.909
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 16256.909
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 48770.727
$ws.Range("N77").Value = -58130.727
